$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the uncertainty/correction columns: the "syst_u" column is dropped
# from the header sequence, "syst1_c" shifts into H, "syst2_c" shifts into I,
# and a new "norm_c" column label takes J. The underlying numeric data in
# H:J is unchanged - only the header text (and hence which shared-string each
# row's K:N labels resolve to) shifts.
$ws.Range("H1").Value = "syst1_c"
$ws.Range("I1").Value = "syst2_c"
$ws.Range("J1").Value = "norm_c"

# Update the sheet's last-known selection/active cell.
$ws.Range("N21").Select()
